$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Total" column header in G1, matching the bold/centered style
# used by the other header cells in row 1
$ws.Range("G1").Value = "Total"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108

# Fill in the Total column values for each realm row (2-6)
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 12
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 17
$ws.Range("G6").Value = 8
